$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("grants por usuario")

# Determine the last used row in column A (data rows of the grants-per-user table)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# Replace every user name in column B with "alice.gomes" (quotes included,
# matching the literal text already used for the other quoted usernames).
$range = $ws.Range($ws.Cells.Item(1, 2), $ws.Cells.Item($lastRow, 2))
$range.Value = '"alice.gomes"'

$wb.Save()
